$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format phone number column (A2:A6) as Text so leading zeros are preserved
$ws.Range("A2:A6").NumberFormat = "@"

# Update the first farmer's phone number
$ws.Range("A2").Value = "05069468693"

# Clear the PackageTier values (column D) that are no longer used
$ws.Range("D2").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("D6").ClearContents()

# Move the active selection to D5
$ws.Range("D5").Select()

$wb.Save()
